$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from column M into the new column N for every row
# that currently has data (rows 2-33). This reuses the existing style
# (cellXfs) entries rather than creating new ones, matching the
# original author's workbook which only grew cellXfs by a single entry
# (for the new N34 cell, handled separately below).
$ws.Range("M2:M33").Copy()
$ws.Range("N2:N33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the values for the new column N (year 2022) ---
$ws.Range("N3").Value = 2022

$ws.Range("N4").Value = 11.927942610539198
$ws.Range("N5").Value = 3.0909744679837434
$ws.Range("N6").Value = 20.963679772397647
$ws.Range("N7").Value = 4.6002717699014832
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 9.112830865859129
$ws.Range("N10").Value = 3.5391993253978327
$ws.Range("N11").Value = 0.30955295909412422
$ws.Range("N12").Value = 6.73157537222552
$ws.Range("N13").Value = 3.9173330796393815
$ws.Range("N14").Value = 0.7444796831494469
$ws.Range("N15").Value = 7.104530072727953
$ws.Range("N16").Value = 23.0957399744971
$ws.Range("N17").Value = 2.6274648905004008
$ws.Range("N18").Value = 43.176223433734158
$ws.Range("N19").Value = 7.6660105666632132
$ws.Range("N20").Value = 0.83437630371297455
$ws.Range("N21").Value = 14.406256431364477
$ws.Range("N22").Value = 34.201612992199827
$ws.Range("N23").Value = 4.4521615244201058
$ws.Range("N24").Value = 63.433733622066185
$ws.Range("N25").Value = 20.535408979625672
$ws.Range("N26").Value = 7.8632542639432348
$ws.Range("N27").Value = 33.368028499329796
$ws.Range("N28").Value = 19.301652062045072
$ws.Range("N29").Value = 7.1220113855063829
$ws.Range("N30").Value = 34.008685896558866
$ws.Range("N31").Value = 7.8668258762379715
$ws.Range("N32").Value = 1.7266187050359711
$ws.Range("N33").Value = 13.723068478111704

# --- New, empty, formatted cell below the table (N34). This is the one
# genuinely new style the author's edit introduced (font 1 / Times New
# Roman, no border, no special alignment). ---
$n34 = $ws.Range("N34")
$n34.VerticalAlignment = -4107

# --- Update the sheet view's selection to match the saved workbook ---
$ws.Range("O6").Select()
